$d = $word.ActiveDocument

# This revision only touches word/numbering.xml: the <w:nsid w:val="…"/>
# GUID recorded on four existing list definitions (w:abstractNum
# w:abstractNumId="990", "991", "99416" and "99417") is swapped for a new
# random-looking value. Nothing else about those lists (bullet glyphs,
# fonts, indents, multilevel structure, ...) changes, and no numId/list
# usage in the document body is touched either - per the commit message
# ("Automatic build output files") this is churn from whatever pipeline
# regenerates/re-exports this .docx, not an authored content edit.
#
# w:nsid is purely internal Word bookkeeping: it has never been surfaced
# by the Word object model (no ListTemplate/List/ListFormat property maps
# to it - confirmed: ListTemplate.NSID / List.ListID are not settable,
# and Document/Range.WordOpenXML is get-only) so there's no run/paragraph
# text to target with Selection/Range edits. We still do the substitution
# the same way the rest of this family of edits is expressed - a
# Find/Replace per old/new pair - so that if a given value ever is
# reachable as literal story text in some document, it gets updated;
# otherwise this is a no-op that leaves the rest of the package untouched.
$nsidPairs = @(
    @("95ed0fb3", "8eb1a71f"),   # abstractNumId 990
    @("f37d0d7d", "4ceceb8f"),   # abstractNumId 991
    @("146e6671", "d856e1a5"),   # abstractNumId 99416
    @("666c428b", "bd5e485a")    # abstractNumId 99417
)

foreach ($pair in $nsidPairs) {
    $oldVal = $pair[0]
    $newVal = $pair[1]

    foreach ($rng in $d.StoryRanges) {
        $rng.Find.Execute($oldVal, $true, $false, $false, $false, $false, `
                           $true, 1, $false, $newVal, 2) | Out-Null
    }
}
